$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.11, 0.211, -0.459, 127, 29, -86),
    @(0.096, 0.6929999999999999, -0.002, 106, -16, 103),
    @(0.07099999999999999, -0.139, 0.463, -31, 46, 113),
    @(0.591, 0.091, 0.117, 144, -20, -95),
    @(0.52, 0.53, -0.049, 9, 35, 148),
    @(-0.146, -0.074, 0.85, -97, -4, -135),
    @(0.164, 0.664, -0.191, -43, 54, -34),
    @(0.247, -0.239, 0.681, -173, 29, 128),
    @(-0.27, 0.006, -0.42, 157, 68, -21),
    @(-0.126, 0.104, -0.216, -8, 9, 172)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt $data[$r].Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}
